# Add `owndat` and `asmt` fields to IC PIN-level reference file output (#841)
#
# The sheet is a two-row "header" export: row 1 holds the source-table name
# (repeated across every column that belongs to it) and row 2 holds the
# field name within that table. We extend both rows with seven new columns:
#   - four more OWNDAT columns: ADDR1, City, State, Zip Code
#   - three new ASMT columns:   VALASM1, VALASM2, VALASM3
#
# New cells must carry the same cell style as the rest of the header
# (fill + thin border) used throughout the sheet, so we clone that format
# from the existing header range before writing the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the existing header formatting (fill/border) onto the new columns
# first, so none of the following Value writes need to invent a new style.
$ws.Range("R1:W2").Copy()
$ws.Range("X1:AD2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 1 — table/source names, left to right.
$ws.Range("X1").Value = "OWNDAT"
$ws.Range("Y1").Value = "OWNDAT"
$ws.Range("Z1").Value = "OWNDAT"
$ws.Range("AA1").Value = "OWNDAT"
$ws.Range("AB1").Value = "ASMT"
$ws.Range("AC1").Value = "ASMT"
$ws.Range("AD1").Value = "ASMT"

# Row 2 — field names, left to right.
$ws.Range("X2").Value = "ADDR1"
$ws.Range("Y2").Value = "City"
$ws.Range("Z2").Value = "State"
$ws.Range("AA2").Value = "Zip Code"
$ws.Range("AB2").Value = "VALASM1"
$ws.Range("AC2").Value = "VALASM2"
$ws.Range("AD2").Value = "VALASM3"
